# Sync attendance_reports - fix "Recorded By" column ordering
# Change "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System"
# across the whole "Session Analysis Results" sheet (column G).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
